# FIX: error in bitbucket conditions
# Slide 3 ("Platforms"), Content Placeholder 2, BitBucket bullet:
# Split the single run
#   ": Free for projects with up to 5 users but private repos are free. Works with Mercurial as well as "
# into:
#   ": Free for projects with up to 5 users but private repos are "
#   "free ("
#   "unlimited for "
#   "academic users)"
#   ". "
#   "Works with Mercurial as well as "

$p = $ppt.ActivePresentation

# Find the slide / shape holding the BitBucket bullet text rather than assuming
# fixed indices, so the script is resilient to minor layout differences.
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($hi = 1; $hi -le $slide.Shapes.Count; $hi++) {
        $shp = $slide.Shapes.Item($hi)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text.IndexOf("Mercurial") -ge 0) {
                $targetShape = $shp
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$fullText = $tr.Text

# The whole original run we need to split (found via plain string search on the
# full shape text, then mapped onto the TextRange via Characters(start,len) --
# this avoids relying on the (unreliable in this host) Runs()/Paragraphs()
# collection Count/Item semantics).
$oldFragment = ": Free for projects with up to 5 users but private repos are free. Works with Mercurial as well as "
$idx0 = $fullText.IndexOf($oldFragment)
$runStart = $idx0 + 1

# Prefix that must stay untouched (keeps the original run / formatting).
$prefix = ": Free for projects with up to 5 users but private repos are "
$prefixLen = $prefix.Length

# Remaining text (after the prefix) that needs to be replaced by the new runs.
$oldTail = $oldFragment.Substring($prefixLen)

# New chunks that replace $oldTail, each becoming its own run.
$chunks = @(
    "free (",
    "unlimited for ",
    "academic users)",
    ". ",
    "Works with Mercurial as well as "
)

# Current absolute (1-based) start of the still-unprocessed remainder of the
# original tail text, and the remainder text itself.
$cursorStart = $runStart + $prefixLen
$remaining = $oldTail

for ($ci = 0; $ci -lt $chunks.Length; $ci++) {
    $chunk = $chunks[$ci]
    if ($ci -lt ($chunks.Length - 1)) {
        # Consume a single old character as the anchor for this split so that
        # the rest of $remaining survives as its own (still original-format) run
        # for the next iteration.
        $consume = 1
    } else {
        # Last chunk absorbs whatever old text is left.
        $consume = $remaining.Length
    }

    $target = $tr.Characters($cursorStart, $consume)
    $target.Text = $chunk

    $remaining = $remaining.Substring($consume)
    $cursorStart = $cursorStart + $chunk.Length
}
